$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 4744.5
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 4744.5
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 4744.5
$ws.Range("M19").ClearContents()
$ws.Range("N19").Value = -5094.5
$ws.Range("H32").Value = 6454.5
$ws.Range("J32").Value = 3074
$ws.Range("L32").Value = 3074
$ws.Range("N32").Value = -3726
$ws.Range("H33").Value = 13586.75
$ws.Range("I33").Value = 16144.3
$ws.Range("K33").Value = 16144.3
$ws.Range("M33").Value = -15915.3
$ws.Range("H39").Value = 312.92856
$ws.Range("I39").Value = 330.5
$ws.Range("J39").Value = 289.5
$ws.Range("K39").Value = 991.5
$ws.Range("L39").Value = 868.5
$ws.Range("M39").Value = -695.5
$ws.Range("N39").Value = -1460.5
$ws.Range("H40").Value = 7030.846
$ws.Range("J40").Value = 8379.556
$ws.Range("L40").Value = 8379.556
$ws.Range("N40").Value = -8729.556
$ws.Range("H116").Value = 4150
$ws.Range("J116").Value = 4720
$ws.Range("L116").Value = 4720
$ws.Range("N116").Value = -11604
$ws.Range("H125").Value = 8123.2104
$ws.Range("I125").Value = 7961.9165
$ws.Range("J125").Value = 8399.714
$ws.Range("K125").Value = 71657.2485
$ws.Range("L125").Value = 75597.42600000001
$ws.Range("M125").Value = -69197.2485
$ws.Range("N125").Value = -80517.42600000001
$ws.Range("H135").Value = 875.95654
$ws.Range("I135").Value = 752.2857
$ws.Range("K135").Value = 6770.571300000001
$ws.Range("M135").Value = -4235.571300000001
$ws.Range("H137").Value = 4395.0454
$ws.Range("I137").Value = 3612
$ws.Range("K137").Value = 10836
$ws.Range("M137").Value = -8286
$ws.Range("H138").Value = 5083.5454
$ws.Range("I138").Value = 2197.111
$ws.Range("J138").Value = 7081.846
$ws.Range("K138").Value = 6591.333
$ws.Range("L138").Value = 21245.538
$ws.Range("M138").Value = -1451.333
$ws.Range("N138").Value = -31525.538
$ws.Range("H141").Value = 3071.923
$ws.Range("I141").Value = 2767.1765
$ws.Range("K141").Value = 8301.529500000001
$ws.Range("M141").Value = -3121.529500000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6183.9824
$ws.Range("I32").Value = 5104.846
$ws.Range("K32").Value = 5104.846
$ws.Range("M32").Value = -4817.846
$ws.Range("H43").Value = 213436.44
$ws.Range("I43").Value = 53655.332
$ws.Range("J43").Value = 293327
$ws.Range("K43").Value = 53655.332
$ws.Range("L43").Value = 293327
$ws.Range("M43").Value = -53342.332
$ws.Range("N43").Value = -293953
$ws.Range("H45").Value = 147599.64
$ws.Range("J45").Value = 6142.857
$ws.Range("L45").Value = 6142.857
$ws.Range("N45").Value = -6896.857
$ws.Range("H74").Value = 3696.0645
$ws.Range("I74").Value = 2079.9614
$ws.Range("K74").Value = 2079.9614
$ws.Range("M74").Value = -1205.9614
$ws.Range("H77").Value = 3696.0645
$ws.Range("I77").Value = 2079.9614
$ws.Range("K77").Value = 10399.807
$ws.Range("M77").Value = -6031.807000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 2528.2144
$ws.Range("I22").Value = 2866.25
$ws.Range("K22").Value = 2866.25
$ws.Range("M22").Value = -2693.25
$ws.Range("H92").Value = 40000
$ws.Range("I92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("M92").ClearContents()
$ws.Range("H134").Value = 5367.839
$ws.Range("I134").Value = 5367.839
$ws.Range("K134").Value = 16103.517
$ws.Range("M134").Value = -13568.517

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 1699.8
$ws.Range("I2").Value = 874.75
$ws.Range("J2").Value = 5000
$ws.Range("K2").Value = 874.75
$ws.Range("L2").Value = 5000
$ws.Range("M2").Value = -761.75
$ws.Range("N2").Value = -5226
$ws.Range("H7").Value = 41.5
$ws.Range("I7").Value = 26.071428
$ws.Range("K7").Value = 26.071428
$ws.Range("M7").Value = 86.928572
$ws.Range("H11").Value = 3158.5715
$ws.Range("I11").Value = 6798.3335
$ws.Range("J11").Value = 428.75
$ws.Range("K11").Value = 6798.3335
$ws.Range("L11").Value = 428.75
$ws.Range("M11").Value = -6658.3335
$ws.Range("N11").Value = -708.75
$ws.Range("H12").Value = 8751.5
$ws.Range("I12").Value = 8500
$ws.Range("J12").Value = 9003
$ws.Range("K12").Value = 8500
$ws.Range("L12").Value = 9003
$ws.Range("M12").Value = -8330
$ws.Range("N12").Value = -9343
$ws.Range("H13").Value = 5000
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 5000
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 5000
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value = -5278

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H48").Value = 600
$ws.Range("I48").Value = 600
$ws.Range("K48").Value = 1800
$ws.Range("M48").Value = -1550
$ws.Range("H68").Value = 2435.2856
$ws.Range("I68").Value = 3118.5
$ws.Range("J68").Value = 2014.8462
$ws.Range("K68").Value = 9355.5
$ws.Range("L68").Value = 6044.5386
$ws.Range("M68").Value = -8544.5
$ws.Range("N68").Value = -7666.5386
$ws.Range("H71").Value = 2435.2856
$ws.Range("I71").Value = 3118.5
$ws.Range("J71").Value = 2014.8462
$ws.Range("K71").Value = 28066.5
$ws.Range("L71").Value = 18133.6158
$ws.Range("M71").Value = -24010.5
$ws.Range("N71").Value = -26245.6158
$ws.Range("H88").Value = 15000
$ws.Range("J88").Value = 15000
$ws.Range("L88").Value = 45000
$ws.Range("N88").Value = -45856
$ws.Range("H91").Value = 15000
$ws.Range("J91").Value = 15000
$ws.Range("L91").Value = 45000
$ws.Range("N91").Value = -47964
$ws.Range("H131").Value = 20835800
$ws.Range("J131").Value = 2698.762
$ws.Range("L131").Value = 8096.286
$ws.Range("N131").Value = -18176.286

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2092.75
$ws.Range("I113").Value = 2174.0908
$ws.Range("J113").Value = 1198
$ws.Range("K113").Value = 2174.0908
$ws.Range("L113").Value = 1198
$ws.Range("M113").Value = -4.090799999999945
$ws.Range("N113").Value = -5538
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H132").Value = 9214.875
$ws.Range("I132").Value = 9102.857
$ws.Range("K132").Value = 27308.571
$ws.Range("M132").Value = -24778.571

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5091.7856
$ws.Range("I7").Value = 2662.2727
$ws.Range("J7").Value = 14000
$ws.Range("K7").Value = 2662.2727
$ws.Range("L7").Value = 14000
$ws.Range("M7").Value = -2550.2727
$ws.Range("N7").Value = -14224
$ws.Range("H22").Value = 1303
$ws.Range("I22").Value = 1551.1428
$ws.Range("K22").Value = 1551.1428
$ws.Range("M22").Value = -1256.1428
$ws.Range("H27").Value = 1303
$ws.Range("I27").Value = 1551.1428
$ws.Range("K27").Value = 1551.1428
$ws.Range("M27").Value = -1444.1428
$ws.Range("H61").Value = 172322.83
$ws.Range("I61").Value = 257484.5
$ws.Range("K61").Value = 257484.5
$ws.Range("M61").Value = -257282.5
$ws.Range("H113").Value = 172322.83
$ws.Range("I113").Value = 257484.5
$ws.Range("K113").Value = 257484.5
$ws.Range("M113").Value = -255314.5
$ws.Range("H126").Value = 5091.7856
$ws.Range("I126").Value = 2662.2727
$ws.Range("J126").Value = 14000
$ws.Range("K126").Value = 7986.8181
$ws.Range("L126").Value = 42000
$ws.Range("M126").Value = -5516.8181
$ws.Range("N126").Value = -46940

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 16395.924
$ws.Range("I45").Value = 7222
$ws.Range("J45").Value = 20473.223
$ws.Range("K45").Value = 7222
$ws.Range("L45").Value = 20473.223
$ws.Range("M45").Value = -6731
$ws.Range("N45").Value = -21455.223
$ws.Range("H107").Value = 431.77777
$ws.Range("I107").Value = 334.1
$ws.Range("J107").Value = 710.8570999999999
$ws.Range("K107").Value = 1002.3
$ws.Range("L107").Value = 2132.5713
$ws.Range("M107").Value = 917.6999999999999
$ws.Range("N107").Value = -5972.5713
$ws.Range("H113").Value = 529.96
$ws.Range("I113").Value = 444.57895
$ws.Range("J113").Value = 800.3333
$ws.Range("K113").Value = 1333.73685
$ws.Range("L113").Value = 2400.9999
$ws.Range("M113").Value = 836.26315
$ws.Range("N113").Value = -6740.9999
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()
$ws.Range("H136").Value = 1791.1666
$ws.Range("I136").Value = 1379.2885
$ws.Range("J136").Value = 12500
$ws.Range("K136").Value = 4137.8655
$ws.Range("L136").Value = 37500
$ws.Range("M136").Value = -1587.8655
$ws.Range("N136").Value = -42600
